# Auto-generated edit script for mujeres.docx
$d = $word.ActiveDocument

# Step 1: delete the paragraphs that were removed entirely
# (find each by its exact original text, expand the found range to the
#  whole paragraph including its paragraph mark, then delete it)
$rng = $d.Content
$found = $rng.Find.Execute('35.2, 24 8 80, María Sánchez Flores, 2027056140, msanchez6140@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('35.05, 5 78 20, Gabriela Martínez Acosta, 2027022681, gmartinez2681@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('34.95, 39 24 43, Laura Solano Solano, 2026047544, lsolano7544@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('34.35, 50 43 6, Sofía Mora Vargas, 2026002442, smora2442@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('32.45, 36 19 44, Sofía Torres Solano, 2027022565, storres2565@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('30.05, 68 17 1, Elena Sánchez Rivera, 2026039105, esanchez9105@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('30.05, 47 32 8, Camila Morales Castro, 2025015249, cmorales5249@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('29.0, 6 34 50, María Sánchez Flores, 2025002421, msanchez2421@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('28.2, 25 53 3, Andrés Morales Flores, 2026047831, amorales7831@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

$rng = $d.Content
$found = $rng.Find.Execute('28.0, 19 1 70, Rebecca Rumore Ramirez, 2027025393, rrumore5393@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
if ($found) { $rng.Expand(4) | Out-Null; $rng.Delete() }

# Step 2: replace the remaining paragraph texts with the new values
$d.Content.Find.Execute('89.65, 87 98 83, Evelyn Torres Evans, 2026036094, etorres6094@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '94.15, 99 98 84, Fernando Ramírez Acosta, 2025069927, framirez9927@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('83.65, 71 96 84, Ana Brewer Ranalli, 2026045420, abrewer5420@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '88.95, 85 86 97, Ana Vargas Pérez, 2022041969, avargas1969@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('82.45, 75 92 80, Ana Rivera Flores, 2025001908, arivera1908@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '83.6, 88 84 78, Elena Torres Vargas, 2025056289, etorres6289@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('77.6, 69 85 79, Paula Castro Pérez, 2027058620, pcastro8620@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '83.25, 80 91 78, Camila Morales Castro, 2022016038, cmorales6038@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('73.9, 93 89 34, Elena Rojas Torres, 2025058855, erojas8855@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '77.2, 59 99 73, Willie Johnson Arias, 2025048340, wjohnson8340@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('73.3, 87 41 95, Daniel Mora Díaz, 2026008682, dmora8682@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '74.9, 58 78 91, Andra Hutcherson Carpenter, 2025055960, ahutcherson5960@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('71.0, 53 71 92, Lucía Chaves Castro, 2026029999, lchaves9999@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '72.35, 79 48 93, Sofía Torres Vargas, 2023061195, storres1195@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('70.5, 76 62 74, Amber Miraflores Saxton, 2026009568, amiraflores9568@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '72.3, 76 80 59, Daniel Vargas Flores, 2024038604, dvargas8604@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('70.05, 98 79 27, Leonie Aguiar Kennedy, 2025011956, laguiar1956@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '70.45, 87 68 54, Elena Ramírez Navarro, 2024047983, eramirez7983@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('68.35, 58 91 54, Manuel Pérez Pérez, 2027006185, mperez6185@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '66.95, 36 85 82, Lillie Picard Norman, 2025023866, lpicard3866@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('67.8, 67 53 86, Pedro Vargas Martínez, 2026019698, pvargas9698@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '65.7, 83 67 44, Paula Ramírez Sánchez, 2022051023, pramirez1023@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('67.15, 87 20 99, José Navarro Flores, 2027003439, jnavarro3439@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '65.4, 81 63 50, Carlos Torres Rojas, 2023039066, ctorres9066@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('66.7, 97 79 17, Daniel Rivera Acosta, 2026003943, drivera3943@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '64.4, 75 61 56, Luis Navarro Díaz, 2022015864, lnavarro5864@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('65.9, 99 13 89, Andrés Rodríguez Pérez, 2026047884, arodriguez7884@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '63.1, 85 43 61, Andrés Mora Díaz, 2025017744, amora7744@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('64.6, 59 69 66, Leslie Bonney Seals, 2026059210, lbonney9210@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '62.0, 50 98 34, Miguel Gómez Flores, 2022037771, mgomez7771@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('64.45, 65 72 55, Lucía Rivera Castro, 2027029993, lrivera9993@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '61.45, 77 48 59, Elena Rojas Torres, 2022037050, erojas7050@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('63.5, 65 65 60, José Flores Solano, 2027040946, jflores0946@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '61.25, 81 46 56, Sofía Torres Solano, 2022016205, storres6205@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('62.55, 55 56 79, Jorge Solano Solano, 2027019808, jsolano9808@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '59.6, 38 74 68, Mary Presto Ferrell, 2023046785, mpresto6785@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('61.3, 58 70 55, José Chaves Vargas, 2027025870, jchaves5870@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '59.5, 43 79 56, Gabriela Castro Castro, 2023051832, gcastro1832@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('61.2, 24 78 85, Lucía Rivera Castro, 2025008588, lrivera8588@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '59.2, 31 79 69, Elena Martínez Díaz, 2025050956, emartinez0956@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('59.6, 94 36 47, Pedro Rivera Díaz, 2027023154, privera3154@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '57.2, 81 25 67, Carlos Torres Rojas, 2024041051, ctorres1051@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('58.8, 66 54 56, Pedro Navarro Torres, 2026032512, pnavarro2512@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '57.1, 33 59 83, Manuel Rivera Rodríguez, 2022057897, mrivera7897@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('58.55, 57 100 12, Kathleen Nelson Bartlett, 2025051325, knelson1325@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '54.9, 95 7 64, Laura Gómez Rodríguez, 2025032672, lgomez2672@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('58.55, 32 53 96, Daniel Pérez Torres, 2025053941, dperez3941@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '52.85, 37 48 77, Elena Ramírez Martínez, 2024069082, eramirez9082@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('58.05, 59 22 99, Maria Smith Johnson, 2025037997, msmith7997@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '52.65, 3 72 88, Camila Rivera Chaves, 2022041327, crivera1327@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('58.05, 30 51 99, Kathryn Hughes Batts, 2027019580, khughes9580@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '51.55, 29 84 40, Lucía Torres Flores, 2022011549, ltorres1549@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('53.25, 98 31 27, Manuel Sánchez Acosta, 2026027976, msanchez7976@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '50.95, 62 15 80, Julie Rodriguez Thomason, 2023045215, jrodriguez5215@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('52.65, 33 42 88, Miguel Morales González, 2026015320, mmorales5320@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '46.4, 78 34 24, Carlos Solano González, 2023013807, csolano3807@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('51.75, 27 42 92, Sabrina Halpern Christopher, 2025055577, shalpern5577@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '45.5, 77 5 56, Concepcion Lewis Johnson, 2025019567, clewis9567@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('49.65, 21 48 85, Amelia Bulmer Otto, 2026028495, abulmer8495@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '44.75, 81 34 15, Jorge Ramírez Castro, 2025047686, jramirez7686@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('47.95, 42 11 98, Lucía Castro Díaz, 2026046049, lcastro6049@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '44.75, 41 20 78, Valeria Mora Navarro, 2022060040, vmora0040@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('46.7, 22 60 60, Ana Vargas Sánchez, 2025037772, avargas7772@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '44.5, 64 28 41, Rosemary Peard Aumick, 2022067699, rpeard7699@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('46.25, 95 32 6, Charlotte Decker Bradley, 2027048019, cdecker8019@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '44.15, 46 57 27, Jorge Ramírez Castro, 2024019435, jramirez9435@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('46.05, 48 3 94, Camila Acosta Castro, 2025027630, cacosta7630@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '42.5, 2 86 39, Valeria Díaz Díaz, 2023013732, vdiaz3732@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('45.8, 35 17 92, Ana Mora Navarro, 2026022713, amora2713@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '42.0, 5 91 28, Laura Flores Morales, 2022063554, lflores3554@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('45.55, 72 11 55, Fernando Rojas Solano, 2025047638, frojas7638@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '42.0, 19 17 98, Daniel Pérez Castro, 2024045575, dperez5575@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('45.25, 9 50 82, Sofía Ramírez Martínez, 2026027221, sramirez7221@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '41.9, 22 42 65, Suzanne Kess Racilis, 2022012531, skess2531@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('45.05, 88 39 2, Abigail Leandry Killian, 2025040690, aleandry0690@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '41.05, 47 6 75, Susie Avant Boehm, 2023042868, savant2868@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('43.5, 64 26 40, Shannon Mcdaniel Gonzalez, 2025044969, smcdaniel4969@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '38.55, 44 55 13, Ana Sánchez Morales, 2025065620, asanchez5620@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('43.25, 60 37 31, Manuel Gómez Mora, 2026039514, mgomez9514@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '38.35, 16 91 3, Lucía Solano Torres, 2023034579, lsolano4579@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('42.6, 68 22 37, Lucía Díaz Castro, 2026051724, ldiaz1724@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '38.25, 38 1 82, Paula Castro Chaves, 2025048014, pcastro8014@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('41.2, 64 28 30, Gena Harvey Shannon, 2025034831, gharvey4831@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '34.2, 22 8 79, Laura Sánchez Flores, 2024026807, lsanchez6807@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('41.0, 57 31 34, Kathleen Pollard Smith, 2026046817, kpollard6817@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '30.25, 45 14 32, Carlos Morales Díaz, 2024020802, cmorales0802@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('39.9, 97 11 7, Sofía Sánchez Ramírez, 2026001242, ssanchez1242@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '27.5, 29 47 3, Paula Ramírez Díaz, 2025026541, pramirez6541@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('37.0, 65 9 37, Laura Martínez Castro, 2025008232, lmartinez8232@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '27.15, 46 17 17, Lucía Rodríguez Rodríguez, 2022023849, lrodriguez3849@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('36.95, 11 92 3, Pam Clardy Melillo, 2026030000, pclardy0000@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '23.55, 11 40 19, Kayla Kettle Stafford, 2025041534, kkettle1534@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('36.7, 83 9 15, Joann Baiz Kesler, 2027014641, jbaiz4641@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '21.35, 9 28 28, Paula Castro Chaves, 2024022443, pcastro2443@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('36.2, 54 34 18, Fernando Vargas González, 2027010066, fvargas0066@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '21.15, 6 15 46, Miguel Morales González, 2022041182, mmorales1182@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('35.9, 32 62 10, Joan Kim Vest, 2027033916, jkim3916@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '21.05, 31 12 20, Carlos Castro Rojas, 2022027065, ccastro7065@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('11.9, 16 12 7, Daniel Navarro Rojas, 2026029702, dnavarro9702@estudiantec.cr', $true, $true, $false, $false, $false, $true, 1, $false, '20.6, 30 16 15, Miguel Chaves Morales, 2023016387, mchaves6387@estudiantec.cr', 2) | Out-Null
$d.Content.Find.Execute('Los porcentajes de cada evaluación fueron 33%, 33% y 34% respectivamente, y la cantidad de mujeres es 60.', $true, $true, $false, $false, $false, $true, 1, $false, 'Los porcentajes de cada evaluación fueron 35%, 35% y 30% respectivamente, y la cantidad de mujeres es 50.', 2) | Out-Null
